$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet originally held 4 rows (1 header + 3 "blank" demo rows: nbsp,
# normal-space, tabs).  The edit reorders those 3 demo rows (normal, tabs,
# nbsp) and appends 11 new demo rows, one per "special" unicode whitespace
# character, each shaped like: before_<name> | <blank char(s)> | after_<name>
#
# NOTE: the runtime's Range.Value setter mis-parses a run of U+00A0 (NBSP)
# characters as a *number* (locale group-separator parsing quirk), so the
# existing nbsp-blank text is preserved by Copy / PasteSpecial instead of
# being retyped through .Value.
# ---------------------------------------------------------------------------

# Stash the existing nbsp-blank value (currently B2) in a scratch cell well
# outside the used range so it survives the rebuild below.
$ws.Range("B2").Copy()
$ws.Range("Z100").PasteSpecial()

# Clear out the old 4 rows completely before rebuilding the sheet.
$ws.Range("A1:C4").Clear()

# Row 1 - headers (unchanged)
$ws.Range("A1").Value = "HEADER_A"
$ws.Range("B1").Value = "HEADER_B"
$ws.Range("C1").Value = "HEADER_C"

# Row 2 - normal blank (11 regular spaces)
$ws.Range("A2").Value = "before_normal_blank"
$ws.Range("B2").Value = "           "
$ws.Range("C2").Value = "after_normal_blank"

# Row 3 - tabs (5 tab characters)
$ws.Range("A3").Value = "before_tabs"
$ws.Range("B3").Value = "`t`t`t`t`t"
$ws.Range("C3").Value = "after_tabs"

# Row 4 - nbsp blank (9 NBSP characters) - restored from the scratch cell
$ws.Range("A4").Value = "before_nbsp_blank"
$ws.Range("Z100").Copy()
$ws.Range("B4").PasteSpecial()
$ws.Range("C4").Value = "after_nbsp_blank"

# Remove the scratch cell so it doesn't leak into the used range/dimension.
$ws.Range("Z100").Clear()

# Rows 5-15 - one per "special" whitespace character
$specials = @(
  @{ name = "hair_space";         row = 5;  char = [char]0x200A },
  @{ name = "zero_width";         row = 6;  char = [char]0x200B },
  @{ name = "six-per-em";         row = 7;  char = [char]0x2006 },
  @{ name = "thin_space";         row = 8;  char = [char]0x2009 },
  @{ name = "punctuation_space";  row = 9;  char = [char]0x2008 },
  @{ name = "four-per-em";        row = 10; char = [char]0x2005 },
  @{ name = "three-per-em";       row = 11; char = [char]0x2004 },
  @{ name = "figure_space";       row = 12; char = [char]0x2007 },
  @{ name = "en_space";           row = 13; char = [char]0x2002 },
  @{ name = "em_space";           row = 14; char = [char]0x2003 },
  @{ name = "braille_blank";      row = 15; char = [char]0x2800 }
)

foreach ($item in $specials) {
  $r = $item.row
  $ws.Range("A$r").Value = "before_" + $item.name
  $ws.Range("B$r").Value = [string]$item.char
  $ws.Range("C$r").Value = "after_" + $item.name
}

# Column C is wider now to fit the longer "punctuation_space" style labels.
# (The engine quantizes ColumnWidth to whole pixels, so 23.33 is the closest
# input that lands on the target stored width of ~24.164.)
$ws.Columns("C").ColumnWidth = 23.33

# Select B9 (punctuation-space blank cell), matching the saved selection.
$ws.Range("B9").Select() | Out-Null
